$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.911942
$ws.Range("H2").Value = 5.735825999999999
$ws.Range("I2").Value = 0.2156379149120961
$ws.Range("J2").Value = 0.2156379149120961
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 102.8289443333334
$ws.Range("N2").Value = 308.486833
$ws.Range("O2").Value = 0.5559120396302444
$ws.Range("P2").Value = 0.5559120396302443
$ws.Range("Q2").Value = 196.602977486562
$ws.Range("R2").Value = 1769.426797379058
$ws.Range("S2").Value = 0.1198757131003964
$ws.Range("T2").Value = 0.1198757131003964

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.911942
$ws.Range("H3").Value = 5.735825999999999
$ws.Range("I3").Value = 0.2156379149120961
$ws.Range("J3").Value = 0.2156379149120961
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3441717873742006
$ws.Range("P3").Value = 0.3441717873742006
$ws.Range("Q3").Value = 121.71925294089
$ws.Range("R3").Value = 1095.47327646801
$ws.Range("S3").Value = 0.07421648660094189
$ws.Range("T3").Value = 0.0742164866009419

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.911942
$ws.Range("H4").Value = 5.735825999999999
$ws.Range("I4").Value = 0.2156379149120961
$ws.Range("J4").Value = 0.2156379149120961
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 18.481835
$ws.Range("N4").Value = 55.445505
$ws.Range("O4").Value = 0.09991617299555507
$ws.Range("P4").Value = 0.09991617299555505
$ws.Range("Q4").Value = 35.33619657357
$ws.Range("R4").Value = 318.02576916213
$ws.Range("S4").Value = 0.02154571521075778
$ws.Range("T4").Value = 0.02154571521075778

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.898253666666667
$ws.Range("H5").Value = 8.694761
$ws.Range("I5").Value = 0.3268788371019294
$ws.Range("J5").Value = 0.3268788371019295
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 102.8289443333334
$ws.Range("N5").Value = 308.486833
$ws.Range("O5").Value = 0.5559120396302444
$ws.Range("P5").Value = 0.5559120396302443
$ws.Range("Q5").Value = 298.024364953546
$ws.Range("R5").Value = 2682.219284581913
$ws.Range("S5").Value = 0.181715881045296
$ws.Range("T5").Value = 0.181715881045296

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.898253666666667
$ws.Range("H6").Value = 8.694761
$ws.Range("I6").Value = 0.3268788371019294
$ws.Range("J6").Value = 0.3268788371019295
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3441717873742006
$ws.Range("P6").Value = 0.3441717873742006
$ws.Range("Q6").Value = 184.5104459967206
$ws.Range("R6").Value = 1660.594013970485
$ws.Range("S6").Value = 0.1125024736201712
$ws.Range("T6").Value = 0.1125024736201712

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.898253666666667
$ws.Range("H7").Value = 8.694761
$ws.Range("I7").Value = 0.3268788371019294
$ws.Range("J7").Value = 0.3268788371019295
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 18.481835
$ws.Range("N7").Value = 55.445505
$ws.Range("O7").Value = 0.09991617299555507
$ws.Range("P7").Value = 0.09991617299555505
$ws.Range("Q7").Value = 53.56504605547833
$ws.Range("R7").Value = 482.085414499305
$ws.Range("S7").Value = 0.03266048243646225
$ws.Range("T7").Value = 0.03266048243646225

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.056250666666667
$ws.Range("H8").Value = 12.168752
$ws.Range("I8").Value = 0.4574832479859744
$ws.Range("J8").Value = 0.4574832479859745
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 102.8289443333334
$ws.Range("N8").Value = 308.486833
$ws.Range("O8").Value = 0.5559120396302444
$ws.Range("P8").Value = 0.5559120396302443
$ws.Range("Q8").Value = 417.099974004713
$ws.Range("R8").Value = 3753.899766042417
$ws.Range("S8").Value = 0.2543204454845519
$ws.Range("T8").Value = 0.2543204454845519

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.056250666666667
$ws.Range("H9").Value = 12.168752
$ws.Range("I9").Value = 0.4574832479859744
$ws.Range("J9").Value = 0.4574832479859745
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 63.66262833333334
$ws.Range("N9").Value = 190.987885
$ws.Range("O9").Value = 0.3441717873742006
$ws.Range("P9").Value = 0.3441717873742006
$ws.Range("Q9").Value = 258.2315786188356
$ws.Range("R9").Value = 2324.08420756952
$ws.Range("S9").Value = 0.1574528271530875
$ws.Range("T9").Value = 0.1574528271530875

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.056250666666667
$ws.Range("H10").Value = 12.168752
$ws.Range("I10").Value = 0.4574832479859744
$ws.Range("J10").Value = 0.4574832479859745
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 18.481835
$ws.Range("N10").Value = 55.445505
$ws.Range("O10").Value = 0.09991617299555507
$ws.Range("P10").Value = 0.09991617299555505
$ws.Range("Q10").Value = 74.96695553997334
$ws.Range("R10").Value = 674.7025998597601
$ws.Range("S10").Value = 0.04570997534833503
$ws.Range("T10").Value = 0.04570997534833503

Write-Output "Done"